$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -------------------------------------------------
# Drop the trailing "Bibliografia" paragraph row entirely (was row 22).
$ws.Rows.Item(22).Delete()

# Row 17 loses its custom 120pt height (back to sheet default 15pt).
# Excel has no "clear custom height" verb, so delete + re-insert a fresh
# blank row in its place, which carries no explicit RowHeight.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()

# Remaining custom row heights that changed.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- Cell content ---------------------------------------------------------
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOB1259"
$ws.Range("C2").Value = "LOB1259"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Introdução a Ecotoxicologia Aquática"
$ws.Range("C3").Value = " Introdução a Ecotoxicologia Aquática"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Introduction to Aquatic Ecotoxicology"
$ws.Range("C4").Value = "Introduction to Aquatic Ecotoxicology"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2020"
$ws.Range("C8").Value = "01/01/2020"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C10").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Introduce theoretical and practical concepts of Aquatic Ecotoxicology for Environmental Engineering students."
$ws.Range("C11").Value = "Introduce theoretical and practical concepts of Aquatic Ecotoxicology for Environmental Engineering students."
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Provide knowledge of theoretical and practical basics about aquatic ecotoxicology and the techniques used in the laboratory for the cultivation and the assays with standardized test organisms."
$ws.Range("C14").Value = "Provide knowledge of theoretical and practical basics about aquatic ecotoxicology and the techniques used in the laboratory for the cultivation and the assays with standardized test organisms."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "A. Theoretical: 1. Ecotoxicology: Introduction, history, concepts; 2. Introduction of chemicals in the aquatic environment: bioavailability of contaminants, synergistic and antagonistic effects, impacts on aquatic systems; 3. Methods for toxicity tests with aquatic organisms: use of bioindicators; B. practice: 4. Selection, maintenance and cultivation of aquatic organisms: good practice; 5.5. Toxicity tests with aquatic organisms."
$ws.Range("C16").Value = "A. Theoretical: 1. Ecotoxicology: Introduction, history, concepts; 2. Introduction of chemicals in the aquatic environment: bioavailability of contaminants, synergistic and antagonistic effects, impacts on aquatic systems; 3. Methods for toxicity tests with aquatic organisms: use of bioindicators; B. practice: 4. Selection, maintenance and cultivation of aquatic organisms: good practice; 5.5. Toxicity tests with aquatic organisms."
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C18").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."
$ws.Range("C19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Nota final: NF ≥ 5,0"
$ws.Range("C21").Value = "Nota final: NF ≥ 5,0"
